# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'87.703.15"
$ws.Range("E2").Value = "  +9.38%  "
$ws.Range("D3").Value = "'3.314.82"
$ws.Range("E3").Value = "  +4.17%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'219.26"
$ws.Range("E5").Value = "  +5.18%  "
$ws.Range("D6").Value = "'635.16"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("D7").Value = "'0.328"
$ws.Range("E7").Value = "  +19.19%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'0.614"
$ws.Range("E9").Value = "  +4.45%  "
$ws.Range("D10").Value = "'3.311.27"
$ws.Range("E10").Value = "  +4.12%  "
$ws.Range("D11").Value = "'0.606"
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").Value = "'0.0000275"
$ws.Range("E12").Value = "  +6.29%  "
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").Value = "'3.939.32"
$ws.Range("E14").Value = "  +5.06%  "
$ws.Range("D15").Value = "'34.22"
$ws.Range("E15").Value = "  +7.67%  "
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "'87.545.40"
$ws.Range("E17").Value = "  +9.57%  "
$ws.Range("D18").Value = "'3.332.17"
$ws.Range("E18").Value = "  +5.34%  "
$ws.Range("E19").Value = "  +6.81%  "
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").Value = "'448.83"
$ws.Range("E21").Value = "  +3.01%  "
$ws.Range("D22").Value = "'9.07"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  +3.71%  "
$ws.Range("D24").Value = "'7.38"
$ws.Range("E24").Value = "  +7.06%  "
$ws.Range("D25").Value = "'5.36"
$ws.Range("E25").Value = "  +14.97%  "
$ws.Range("D26").Value = "'12.35"
$ws.Range("E26").Value = "  +13.44%  "
$ws.Range("D27").Value = "'3.509.52"
$ws.Range("E27").Value = "  +5.23%  "
$ws.Range("D28").Value = "'78.57"
$ws.Range("E28").Value = "  +3.68%  "
$ws.Range("B29").Value = "Cronos"
$ws.Range("C29").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D29").Value = "'0.215"
$ws.Range("E29").Value = "  +74.82%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0000129"
$ws.Range("E30").Value = "  +5.38%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'596.73"
$ws.Range("E32").Value = "  +7.39%  "
$ws.Range("D33").Value = "'9.29"
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D35").Value = "'1.53"
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("D37").Value = "'0.153"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").Value = "'23.47"
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("E39").Value = "  +17.31%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.417"
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("D42").Value = "'21.40"
$ws.Range("E42").Value = "  +3.02%  "
$ws.Range("D43").Value = "'2.05"
$ws.Range("E43").Value = "  +13.56%  "
$ws.Range("E44").Value = "  +13.42%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'157.20"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("D47").Value = "'188.99"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").Value = "'46.31"
$ws.Range("E48").Value = "  +8.31%  "
$ws.Range("E49").Value = "  +5.53%  "
$ws.Range("D50").Value = "'0.783"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'26.47"
$ws.Range("E51").Value = "  +7.60%  "
